# Add results of FFNN day ahead with WD and time
$wb = $excel.ActiveWorkbook

$wsDay = $wb.Worksheets.Item("Day Ahead")
$wsIntra = $wb.Worksheets.Item("Intra Day")

# --- "Day Ahead" sheet: fill in new RMSE values -------------------------
# Row 2 (FFNN / 1 Forecast WS): C2 = 0.181 (reuse the existing red-font style
# already used for the highlighted cells on the "Intra Day" sheet)
$wsDay.Range("C2").Value = 0.181
$wsIntra.Range("C3").Copy() | Out-Null
$wsDay.Range("C2").PasteSpecial(-4122) | Out-Null

# Row 3 (FFNN / 1 Forecast WS + WD): C3 = 0.177, D3 = 0.179 (+ threaded comment)
$wsDay.Range("C3").Value = 0.177
$wsDay.Range("D3").Value = 0.179

# Row 4 (FFNN / 1 Forecast WS + WD + time): C4 = 0.168
$wsDay.Range("C4").Value = 0.168

# New row 8: LSTM / "24 past power"
$wsDay.Range("A8").Value = "LSTM"
$wsDay.Range("B8").Value = "24 past power"

# Threaded comment on D3
$wsDay.Range("D3").AddCommentThreaded("With u an v for wind direction. Same results again without ws, but only u and v.") | Out-Null

# Selections, matching the saved workbook state
$wsDay.Range("C6").Select() | Out-Null
$wsIntra.Range("G20").Select() | Out-Null
